# Daily attendance processing - 2025-10-28 21:20:09
# Reorders the comma-separated "Recorded By" list in column G so that the
# actual editor (email / service account) is listed first and "System" /
# "system" is moved to the end, e.g.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"         -> "backup@backdoor.com, System"
#   "system, System, backup@backdoor.com" -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")
if ($null -eq $ws) {
    $ws = $wb.ActiveSheet
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "system, System, backup@backdoor.com" = "backup@backdoor.com, System, system";
}

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
